$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.421.63"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "2.375.84"
$ws.Range("E3").Value = "  +5.70%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'235.43"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").Value = "'0.656"
$ws.Range("E6").Value = "  +2.27%  "
$ws.Range("D7").Value = "'72.16"
$ws.Range("E7").Value = "  +14.61%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +4.03%  "
$ws.Range("D10").Value = "'0.0973"
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("D11").Value = "'57.30"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "'26.75"
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("D13").Value = "2.735.64"
$ws.Range("E13").Value = "  +5.88%  "
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").Value = "'15.83"
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("D16").Value = "'6.25"
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("D17").Value = "'0.856"
$ws.Range("E17").Value = "  +3.43%  "
$ws.Range("D18").Value = "2.377.72"
$ws.Range("E18").Value = "  +4.13%  "
$ws.Range("D19").Value = "43.432.29"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("D20").Value = "0.0₃0992"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D21").Value = "'6.39"
$ws.Range("E21").Value = "  +5.50%  "
$ws.Range("D22").Value = "'74.55"
$ws.Range("E22").Value = "  +2.57%  "
$ws.Range("D23").Value = "'251.72"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D24").Value = "'3.95"
$ws.Range("E24").Value = "  +18.55%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "'2.47"
$ws.Range("E26").Value = "  +1.88%  "
$ws.Range("D27").Value = "'23.11"
$ws.Range("E27").Value = "  +10.04%  "
$ws.Range("D28").Value = "'10.03"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("D30").Value = "'174.43"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("E31").Value = "  +8.90%  "
$ws.Range("E32").Value = "  -9.01%  "
$ws.Range("D33").Value = "'0.128"
$ws.Range("E33").Value = "  +3.40%  "
$ws.Range("E34").Value = "  +4.12%  "
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("D36").Value = "'5.09"
$ws.Range("E36").Value = "  +2.97%  "
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").Value = "'6.62"
$ws.Range("E37").Value = "  +3.57%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.46"
$ws.Range("E38").Value = "  +8.37%  "
$ws.Range("D39").Value = "'3.66"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("D41").Value = "'8.94"
$ws.Range("E41").Value = "  +3.53%  "
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("D43").Value = "'18.66"
$ws.Range("E43").Value = "  +9.90%  "
$ws.Range("E44").Value = "  +10.42%  "
$ws.Range("E45").Value = "  +2.48%  "
$ws.Range("E46").Value = "  +2.49%  "
$ws.Range("D47").Value = "'4.48"
$ws.Range("E47").Value = "  +3.13%  "
$ws.Range("D48").Value = "'0.0953"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").Value = "1.454.73"
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("D50").Value = "2.602.81"
$ws.Range("E50").Value = "  +5.96%  "
$ws.Range("E51").Value = "  -0.62%  "
